$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 56.216147
$ws.Range("H2").Value = 168.648441
$ws.Range("I2").Value = 0.9695233148109074
$ws.Range("J2").Value = 0.9695233148109074
$ws.Range("M2").Value = 3.624854
$ws.Range("N2").Value = 10.874562
$ws.Range("O2").Value = 0.2900317783616697
$ws.Range("P2").Value = 0.2900317783616697
$ws.Range("Q2").Value = 203.775325317538
$ws.Range("R2").Value = 1833.977927857842
$ws.Range("S2").Value = 0.2811925711577084
$ws.Range("T2").Value = 0.2811925711577084

# Row 3
$ws.Range("G3").Value = 56.216147
$ws.Range("H3").Value = 168.648441
$ws.Range("I3").Value = 0.9695233148109074
$ws.Range("J3").Value = 0.9695233148109074
$ws.Range("O3").Value = 0.3900977855855255
$ws.Range("P3").Value = 0.3900977855855255
$ws.Range("Q3").Value = 274.0813562306083
$ws.Range("R3").Value = 2466.732206075475
$ws.Range("S3").Value = 0.3782088981812733
$ws.Range("T3").Value = 0.3782088981812733

# Row 4
$ws.Range("G4").Value = 56.216147
$ws.Range("H4").Value = 168.648441
$ws.Range("I4").Value = 0.9695233148109074
$ws.Range("J4").Value = 0.9695233148109074
$ws.Range("M4").Value = 2.893069333333333
$ws.Range("N4").Value = 8.679207999999999
$ws.Range("O4").Value = 0.2314802316645793
$ws.Range("P4").Value = 0.2314802316645793
$ws.Range("Q4").Value = 162.6372109238586
$ws.Range("R4").Value = 1463.734898314728
$ws.Range("S4").Value = 0.2244254815166397
$ws.Range("T4").Value = 0.2244254815166397

# Row 5
$ws.Range("G5").Value = 56.216147
$ws.Range("H5").Value = 168.648441
$ws.Range("I5").Value = 0.9695233148109074
$ws.Range("J5").Value = 0.9695233148109074
$ws.Range("M5").Value = 1.104712
$ws.Range("N5").Value = 3.314136
$ws.Range("O5").Value = 0.08839020438822554
$ws.Range("P5").Value = 0.08839020438822554
$ws.Range("Q5").Value = 62.102652184664
$ws.Range("R5").Value = 558.9238696619759
$ws.Range("S5").Value = 0.08569636395528603
$ws.Range("T5").Value = 0.08569636395528603

# Row 6
$ws.Range("I6").Value = 0.001287520467967504
$ws.Range("J6").Value = 0.001287520467967504
$ws.Range("M6").Value = 3.624854
$ws.Range("N6").Value = 10.874562
$ws.Range("O6").Value = 0.2900317783616697
$ws.Range("P6").Value = 0.2900317783616697
$ws.Range("Q6").Value = 0.2706122670853333
$ws.Range("R6").Value = 2.435510403768
$ws.Range("S6").Value = 0.0003734218510016645
$ws.Range("T6").Value = 0.0003734218510016645

# Row 7
$ws.Range("I7").Value = 0.001287520467967504
$ws.Range("J7").Value = 0.001287520467967504
$ws.Range("O7").Value = 0.3900977855855255
$ws.Range("P7").Value = 0.3900977855855255
$ws.Range("S7").Value = 0.000502258883450163
$ws.Range("T7").Value = 0.000502258883450163

# Row 8
$ws.Range("I8").Value = 0.001287520467967504
$ws.Range("J8").Value = 0.001287520467967504
$ws.Range("M8").Value = 2.893069333333333
$ws.Range("N8").Value = 8.679207999999999
$ws.Range("O8").Value = 0.2314802316645793
$ws.Range("P8").Value = 0.2314802316645793
$ws.Range("Q8").Value = 0.2159811267235555
$ws.Range("R8").Value = 1.943830140512
$ws.Range("S8").Value = 0.0002980355361980054
$ws.Range("T8").Value = 0.0002980355361980055

# Row 9
$ws.Range("I9").Value = 0.001287520467967504
$ws.Range("J9").Value = 0.001287520467967504
$ws.Range("M9").Value = 1.104712
$ws.Range("N9").Value = 3.314136
$ws.Range("O9").Value = 0.08839020438822554
$ws.Range("P9").Value = 0.08839020438822554
$ws.Range("Q9").Value = 0.08247190612266665
$ws.Range("R9").Value = 0.742247155104
$ws.Range("S9").Value = 0.0001138041973176715
$ws.Range("T9").Value = 0.0001138041973176715

# Row 10
$ws.Range("G10").Value = 0.4660483333333333
$ws.Range("H10").Value = 1.398145
$ws.Range("I10").Value = 0.008037632408272877
$ws.Range("J10").Value = 0.008037632408272877
$ws.Range("M10").Value = 3.624854
$ws.Range("N10").Value = 10.874562
$ws.Range("O10").Value = 0.2900317783616697
$ws.Range("P10").Value = 0.2900317783616697
$ws.Range("Q10").Value = 1.689357165276667
$ws.Range("R10").Value = 15.20421448749
$ws.Range("S10").Value = 0.002331168821188773
$ws.Range("T10").Value = 0.002331168821188773

# Row 11
$ws.Range("G11").Value = 0.4660483333333333
$ws.Range("H11").Value = 1.398145
$ws.Range("I11").Value = 0.008037632408272877
$ws.Range("J11").Value = 0.008037632408272877
$ws.Range("O11").Value = 0.3900977855855255
$ws.Range("P11").Value = 0.3900977855855255
$ws.Range("Q11").Value = 2.272214765430555
$ws.Range("R11").Value = 20.449932888875
$ws.Range("S11").Value = 0.003135462603817704
$ws.Range("T11").Value = 0.003135462603817704

# Row 12
$ws.Range("G12").Value = 0.4660483333333333
$ws.Range("H12").Value = 1.398145
$ws.Range("I12").Value = 0.008037632408272877
$ws.Range("J12").Value = 0.008037632408272877
$ws.Range("M12").Value = 2.893069333333333
$ws.Range("N12").Value = 8.679207999999999
$ws.Range("O12").Value = 0.2314802316645793
$ws.Range("P12").Value = 0.2314802316645793
$ws.Range("Q12").Value = 1.348310141017778
$ws.Range("R12").Value = 12.13479126916
$ws.Range("S12").Value = 0.001860553011901736
$ws.Range("T12").Value = 0.001860553011901736

# Row 13
$ws.Range("G13").Value = 0.4660483333333333
$ws.Range("H13").Value = 1.398145
$ws.Range("I13").Value = 0.008037632408272877
$ws.Range("J13").Value = 0.008037632408272877
$ws.Range("M13").Value = 1.104712
$ws.Range("N13").Value = 3.314136
$ws.Range("O13").Value = 0.08839020438822554
$ws.Range("P13").Value = 0.08839020438822554
$ws.Range("Q13").Value = 0.5148491864133333
$ws.Range("R13").Value = 4.63364267772
$ws.Range("S13").Value = 0.000710447971364665
$ws.Range("T13").Value = 0.000710447971364665

# Row 14
$ws.Range("G14").Value = 1.226435333333333
$ws.Range("H14").Value = 3.679306
$ws.Range("I14").Value = 0.02115153231285227
$ws.Range("J14").Value = 0.02115153231285228
$ws.Range("M14").Value = 3.624854
$ws.Range("N14").Value = 10.874562
$ws.Range("O14").Value = 0.2900317783616697
$ws.Range("P14").Value = 0.2900317783616697
$ws.Range("Q14").Value = 4.445649023774666
$ws.Range("R14").Value = 40.01084121397199
$ws.Range("S14").Value = 0.006134616531770865
$ws.Range("T14").Value = 0.006134616531770867

# Row 15
$ws.Range("G15").Value = 1.226435333333333
$ws.Range("H15").Value = 3.679306
$ws.Range("I15").Value = 0.02115153231285227
$ws.Range("J15").Value = 0.02115153231285228
$ws.Range("O15").Value = 0.3900977855855255
$ws.Range("P15").Value = 0.3900977855855255
$ws.Range("Q15").Value = 5.979475247372221
$ws.Range("R15").Value = 53.81527722635
$ws.Range("S15").Value = 0.00825116591698436
$ws.Range("T15").Value = 0.008251165916984361

# Row 16
$ws.Range("G16").Value = 1.226435333333333
$ws.Range("H16").Value = 3.679306
$ws.Range("I16").Value = 0.02115153231285227
$ws.Range("J16").Value = 0.02115153231285228
$ws.Range("M16").Value = 2.893069333333333
$ws.Range("N16").Value = 8.679207999999999
$ws.Range("O16").Value = 0.2314802316645793
$ws.Range("P16").Value = 0.2314802316645793
$ws.Range("Q16").Value = 3.548162452183111
$ws.Range("R16").Value = 31.933462069648
$ws.Range("S16").Value = 0.004896161599839879
$ws.Range("T16").Value = 0.004896161599839879

# Row 17
$ws.Range("G17").Value = 1.226435333333333
$ws.Range("H17").Value = 3.679306
$ws.Range("I17").Value = 0.02115153231285227
$ws.Range("J17").Value = 0.02115153231285228
$ws.Range("M17").Value = 1.104712
$ws.Range("N17").Value = 3.314136
$ws.Range("O17").Value = 0.08839020438822554
$ws.Range("P17").Value = 0.08839020438822554
$ws.Range("Q17").Value = 1.354857829957333
$ws.Range("R17").Value = 12.193720469616
$ws.Range("S17").Value = 0.001869588264257169
$ws.Range("T17").Value = 0.001869588264257169
